# Update Pre outcome measures on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New answer options used for column B (rows 2-15)
$ws.Range("B2").Value  = "Not worse"
$ws.Range("B3").Value  = "Somewhat worse"
$ws.Range("B4").Value  = "Somewhat worse"
$ws.Range("B5").Value  = "A lot worse"
$ws.Range("B6").Value  = "A lot worse"
$ws.Range("B7").Value  = "A lot worse"
$ws.Range("B8").Value  = "A little worse"
$ws.Range("B9").Value  = "A little worse"
$ws.Range("B10").Value = "A little worse"
$ws.Range("B11").Value = "Somewhat worse"
$ws.Range("B12").Value = "Somewhat worse"
$ws.Range("B13").Value = "A little worse"
$ws.Range("B14").Value = "A little worse"
$ws.Range("B15").Value = "A lot worse"
